$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "imageUrl"

$ws.Range("J2").Value = "https://lh3.googleusercontent.com/a/ACg8ocIQfe1aUqsoSuK9TLzuzLVKQAArLvvNnHK2pqi8wQHw7A=s288-c-no"

$ws.Columns.Item(10).ColumnWidth = 101.1666666666667

$excel.CutCopyMode = $false

[void]$ws.Range("B3").Select()
